$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet 1: LP1912
# ---------------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("LP1912")

# Header refresh
$ws1.Range("A2").Value = "Última actualización: 05:53:46"
$ws1.Range("A3").Value = "Total filas: 40"

# A new scrape row is inserted before the existing row 34, pushing the old
# rows 34-37 down to 35-38.
$ws1.Rows.Item(34).Insert()

$ws1.Range("A34").Value = "05:53:46"
$ws1.Range("B34").Value = "07:04"
$ws1.Range("C34").Value = "23_HERNANDEZ"
$ws1.Range("D34").Value = 71
$ws1.Range("E34").Value = "LP1912"

# Seven brand-new rows appended at the bottom of the sheet (39-45).
$ws1.Range("A39").Value = "05:53:46"
$ws1.Range("B39").Value = "07:21"
$ws1.Range("C39").Value = "26_HERNANDEZ"
$ws1.Range("D39").Value = 88
$ws1.Range("E39").Value = "LP1912"

$ws1.Range("A40").Value = "05:53:46"
$ws1.Range("B40").Value = "07:31"
$ws1.Range("C40").Value = "11_ETCHEVERRY"
$ws1.Range("D40").Value = 98
$ws1.Range("E40").Value = "LP1912"

$ws1.Range("A41").Value = "05:53:46"
$ws1.Range("B41").Value = "07:32"
$ws1.Range("C41").Value = "84_COLONIA URQUIZA-ESC 49"
$ws1.Range("D41").Value = 99
$ws1.Range("E41").Value = "LP1912"

$ws1.Range("A42").Value = "05:53:46"
$ws1.Range("B42").Value = "07:36"
$ws1.Range("C42").Value = "27_EL RETIRO"
$ws1.Range("D42").Value = 103
$ws1.Range("E42").Value = "LP1912"

$ws1.Range("A43").Value = "05:53:46"
$ws1.Range("B43").Value = "07:39"
$ws1.Range("C43").Value = "10_OLMOS"
$ws1.Range("D43").Value = 106
$ws1.Range("E43").Value = "LP1912"

$ws1.Range("A44").Value = "05:53:46"
$ws1.Range("B44").Value = "07:47"
$ws1.Range("C44").Value = "14_ABASTO"
$ws1.Range("D44").Value = 114
$ws1.Range("E44").Value = "LP1912"

$ws1.Range("A45").Value = "05:53:46"
$ws1.Range("B45").Value = "07:51"
$ws1.Range("C45").Value = "215D_EL PATO"
$ws1.Range("D45").Value = 118
$ws1.Range("E45").Value = "LP1912"

# ---------------------------------------------------------------------------
# Sheet 2: LP1912-215
# ---------------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("LP1912-215")

$ws2.Range("A2").Value = "Última actualización: 05:53:46"
$ws2.Range("A3").Value = "Total filas: 12"

$ws2.Range("A17").Value = "05:53:46"
$ws2.Range("B17").Value = "07:51"
$ws2.Range("C17").Value = "215D_EL PATO"
$ws2.Range("D17").Value = 118
$ws2.Range("E17").Value = "LP1912"

# ---------------------------------------------------------------------------
# Sheet 3: 6203-6173
# ---------------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("6203-6173")

$ws3.Range("A2").Value = "Última actualización: 05:53:46"
$ws3.Range("A3").Value = "Total filas: 9"

$ws3.Range("A14").Value = "05:53:46"
$ws3.Range("B14").Value = "07:35"
$ws3.Range("C14").Value = "215A_LA PLATA"
$ws3.Range("D14").Value = 102
$ws3.Range("E14").Value = "L6173"
